# "update lgbm for qoq"
#
# A new LGBM result row is inserted into the "average_mae" comparison
# table at row 63 (just above the existing "ibes_1|ni|rnn_double|all"
# row), pushing the two rows that used to be 63-64 down to 64-65. The
# new row carries a fresh label and fresh metric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("average_mae")

# Make room for the new row by shifting the existing rows 63 and 64
# down to 64 and 65 (copy whole rows, bottom-most first, so formatting
# - including the bold/border index-column style - comes along for the
# ride).
$ws.Range("A64:I64").Copy($ws.Range("A65"))
$ws.Range("A63:I63").Copy($ws.Range("A64"))

# Populate the newly freed row 63 with the new comparison result.
$ws.Range("A63").Value = "ibes_1|fwdepsqcut-46|dense2｜compare_hyperopt -code 0 -exclude_fwd True"
$ws.Range("B63").Value = 0.009425532281481011
$ws.Range("C63").Value = 0.008679017950585495
$ws.Range("D63").Value = 0.0001866655144276845
$ws.Range("E63").Value = 0.0001903704215854449
$ws.Range("F63").Value = 0.1436124779217507
$ws.Range("G63").Value = 0.1266150358923953
$ws.Range("H63").Value = 0.258876568546933
$ws.Range("I63").Value = 14156
